# Updated cryptos list - applies price/volume/name changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All these columns (B-E) are stored as plain text in the sheet (coin
# names, links, and numeric-looking price/percentage strings). Force the
# Text number format first so Excel does not auto-convert values such as
# "1.00" or "3.40" into numbers (which would drop trailing zeros).
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '69.391.31'
$ws.Cells.Item(2, 5).Value = '  +1.34%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.947.92'
$ws.Cells.Item(3, 5).Value = '  +0.38%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '503.64'
$ws.Cells.Item(5, 5).Value = '  +3.21%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '148.38'
$ws.Cells.Item(6, 5).Value = '  -0.25%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.627'
$ws.Cells.Item(7, 5).Value = '  -0.49%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.30%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.176'
$ws.Cells.Item(10, 5).Value = '  +3.99%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -1.40%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '43.71'
$ws.Cells.Item(12, 5).Value = '  +1.33%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '10.53'
$ws.Cells.Item(13, 5).Value = '  -1.96%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '4.573.41'
$ws.Cells.Item(14, 5).Value = '  +0.11%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.949.54'
$ws.Cells.Item(15, 5).Value = '  +0.36%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '14.33'
$ws.Cells.Item(16, 5).Value = '  -2.74%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -0.42%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +4.94%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '20.03'
$ws.Cells.Item(19, 5).Value = '  +0.13%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '69.431.99'
$ws.Cells.Item(20, 5).Value = '  +1.29%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '438.46'
$ws.Cells.Item(21, 5).Value = '  -1.62%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -2.05%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '14.66'
$ws.Cells.Item(23, 5).Value = '  -2.62%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '89.02'

# Row 25
$ws.Cells.Item(25, 5).Value = '  +5.92%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '3.87'
$ws.Cells.Item(26, 5).Value = '  +6.18%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '11.23'
$ws.Cells.Item(27, 5).Value = '  -2.27%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '37.18'
$ws.Cells.Item(28, 5).Value = '  -4.50%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '5.66'
$ws.Cells.Item(29, 5).Value = '  -3.42%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '700.72'
$ws.Cells.Item(30, 5).Value = '  -3.71%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '13.48'
$ws.Cells.Item(31, 5).Value = '  -1.73%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.58%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -0.64%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '0.452'
$ws.Cells.Item(34, 5).Value = '  +12.88%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '63.57'
$ws.Cells.Item(35, 5).Value = '  +3.73%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '0.0₃0894'
$ws.Cells.Item(36, 5).Value = '  -1.38%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -3.47%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -2.98%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.151'
$ws.Cells.Item(39, 5).Value = '  +0.71%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 5).Value = '  -0.09%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.02%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.0490'
$ws.Cells.Item(42, 5).Value = '  +1.11%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -1.98%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -3.90%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +2.63%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.94%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +6.77%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '0.0₆0355'
$ws.Cells.Item(48, 5).Value = '  +0.98%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Stacks'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(49, 4).Value = '3.01'
$ws.Cells.Item(49, 5).Value = '  +5.15%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(50, 4).Value = '3.40'
$ws.Cells.Item(50, 5).Value = '  -0.77%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '2.11'
$ws.Cells.Item(51, 5).Value = '  -2.29%  '
